$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 12.65830550497259

$ws.Range("B3").Value = 0.00001292064567892659
$ws.Range("C3").Value = 0.002571899574220771
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 1.249564624054534

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.1494219747398047
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 5.586269137925634

$ws.Range("B5").Value = 0.01293466051926884
$ws.Range("C5").Value = 0.00006240767534437808
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 10.19245300693656
$ws.Range("G5").Value = 10.95819334290503

$ws.Range("B6").Value = 1.455362044514542
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 3.537761648806719
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 7.143138311642302
